$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new rows (pushing existing data down) --------------------
# New row for "zeekr 001 2024" goes right after current row 2 (becomes row 3)
$ws.Rows("3:3").Insert()
# New row for "maxus mifa 7 2024" goes right after that (becomes row 4)
$ws.Rows("4:4").Insert()
# New row for "bmw x2 2022" is inserted later, which (after the two inserts
# above) lands at row 8
$ws.Rows("8:8").Insert()
# New row for "NIO EL6 2024" is inserted further down, landing at row 16
$ws.Rows("16:16").Insert()

# --- Rename car names (hyphens -> spaces) on the rows that already existed
$ws.Range("A2").Value = "ford tourneo custom 2024"
$ws.Range("A5").Value = "ford tourneo custom 2024"
$ws.Range("A6").Value = "vw passat 2024"
$ws.Range("A7").Value = "skoda kodiaq 2024"
$ws.Range("A9").Value = "renault rafale hev 2022"
$ws.Range("A10").Value = "mercedes benz e class 2024"
$ws.Range("A11").Value = "suzuki swift 2024"
$ws.Range("A12").Value = "dacia duster 2024"
$ws.Range("A13").Value = "renault espace 2022"
$ws.Range("A14").Value = "Toyota C HR 2024"
$ws.Range("A15").Value = "Honda CR V 2024"
$ws.Range("A17").Value = "Honda CR V 2024"

# --- Fill in the brand-new rows --------------------------------------------
function Set-DataRow($row, $values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

Set-DataRow 3  @("zeekr 001 2024",   5.95, 27.78, 19.84, 12.3,  13.49, 16.67, 3.97, 0,    100)
Set-DataRow 4  @("maxus mifa 7 2024",0,    23.04, 47.39, 7.83,  9.57,  6.09,  6.09, 0,    100)
Set-DataRow 8  @("bmw x2 2022",      0,    55.81, 20.35, 11.63, 4.65,  5.23,  2.33, 0,    100)
Set-DataRow 16 @("NIO EL6 2024",     0,    38.76, 27.52, 12.02, 5.43,  12.4,  3.88, 0,    100)

# --- Column A width changes from 29 to 28 -----------------------------------
# (The COM ColumnWidth property and the raw OOXML "width" attribute differ by
# a constant offset of 5/6 of a character for the default font, so we adjust
# for that here to land exactly on width="28" in the saved file.)
$ws.Columns("A:A").ColumnWidth = 27.1666666666667
